$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time range note in D11 (append 18:30 end time)
$ws.Range("D11").Value = "13:00-15:00; 16:30-18:30"

# Add the time spent (minutes) and work description for the new entry
$ws.Range("B11").Value = 240
$ws.Range("C11").Value = "Implemented PostUser; made server/dbManager singleton/misc."

# Move the active selection to B12 as in the saved workbook
$ws.Range("B12").Select()

$wb.Save()
